$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.663.80"
$ws.Range("E2").Value = "  -1.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.352.89"
$ws.Range("E3").Value = "  -1.09%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.08"
$ws.Range("E5").Value = "  -1.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.47"
$ws.Range("E6").Value = "  +3.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.635"
$ws.Range("E7").Value = "  -1.21%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -4.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.71"
$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0926"
$ws.Range("E11").Value = "  -1.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.48"
$ws.Range("E12").Value = "  -0.78%  "

$ws.Range("E13").Value = "  -2.20%  "

$ws.Range("E14").Value = "  +0.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.94"
$ws.Range("E15").Value = "  -5.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.710.06"
$ws.Range("E16").Value = "  -0.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.375.52"
$ws.Range("E17").Value = "  -0.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.674.83"
$ws.Range("E18").Value = "  -1.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.72"
$ws.Range("E19").Value = "  -0.57%  "

$ws.Range("E20").Value = "  -1.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.51"
$ws.Range("E21").Value = "  +0.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.58"
$ws.Range("E22").Value = "  +5.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "256.94"
$ws.Range("E23").Value = "  -7.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.33"
$ws.Range("E24").Value = "  -2.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.41"
$ws.Range("E25").Value = "  -1.82%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.45"
$ws.Range("E27").Value = "  -1.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.82"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("E29").Value = "  +1.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.81"
$ws.Range("E30").Value = "  -0.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.71"
$ws.Range("E31").Value = "  -3.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0891"
$ws.Range("E32").Value = "  -3.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.08"
$ws.Range("E33").Value = "  +4.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.90"
$ws.Range("E34").Value = "  -9.49%  "

$ws.Range("E35").Value = "  +20.01%  "

$ws.Range("E36").Value = "  -1.20%  "

$ws.Range("E37").Value = "  -4.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0363"
$ws.Range("E38").Value = "  -0.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.83"
$ws.Range("E39").Value = "  -6.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.67"
$ws.Range("E40").Value = "  -5.04%  "

$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.11"
$ws.Range("E41").Value = "  +4.71%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.238"
$ws.Range("E42").Value = "  +3.20%  "

$ws.Range("E43").Value = "  -6.10%  "

$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "113.75"
$ws.Range("E45").Value = "  -8.27%  "

$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.00"
$ws.Range("E46").Value = "  -3.56%  "

$ws.Range("E47").Value = "  -1.71%  "

$ws.Range("E48").Value = "  -3.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.29"
$ws.Range("E49").Value = "  -14.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.11"
$ws.Range("E50").Value = "  +3.55%  "

$ws.Range("E51").Value = "  -2.41%  "

